$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Internal Assignment" column header in O4, matching the look of the
# other header cells (M4/N4) but bold and with a slightly larger font (size 12).
$ws.Range("O4").Value = "Internal Assignment"
$ws.Range("O4").Font.Bold = $true
$ws.Range("O4").Font.Size = 12
$ws.Range("O4").Font.Name = "Calibri"
$ws.Range("O4").Font.Color = 0

# Populate the new column's data rows (5-7) with the literal text "FALSE" (matching
# the style used by the existing "Mandatory"/"Unique" columns K and L). Typing the
# word FALSE directly into a cell causes Excel to store it as a boolean, so instead
# we build it via a formula and then paste back just the value, which keeps it as
# plain text - exactly like the neighbouring K/L cells.
$ws.Range("Q1").Formula = '="FALSE"'
$ws.Range("Q1").Copy() | Out-Null

$ws.Range("O5").PasteSpecial(-4163)
$ws.Range("O6").PasteSpecial(-4163)
$ws.Range("O7").PasteSpecial(-4163)

# Clean up the scratch cell used to build the text value.
$ws.Range("Q1").Value = ""

# Match the look of the existing data columns (e.g. K, "Mandatory") for these new cells.
$ws.Range("O5:O7").Font.Name = "Calibri"
$ws.Range("O5:O7").Font.Size = 11
$ws.Range("O5:O7").Font.Bold = $false
